$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D from Excel's automatic text-to-number coercion while
# we write values that look numeric (e.g. "1.00", "66.90"); Excel would
# otherwise silently convert these into real numbers losing the literal
# textual representation. We restore the original "Normal" style
# afterwards so no style index changes leak into the saved file.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "46.888.91"
$ws.Range("E2").Value = "  +5.90%  "
$ws.Range("D3").Value = "2.307.59"
$ws.Range("E3").Value = "  +3.88%  "
$ws.Range("E4").Value = "  -0.67%  "
$ws.Range("D5").Value = "303.59"
$ws.Range("E5").Value = "  +1.63%  "
$ws.Range("D6").Value = "101.56"
$ws.Range("E6").Value = "  +13.38%  "
$ws.Range("D7").Value = "0.571"
$ws.Range("E7").Value = "  +2.13%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("E9").Value = "  +8.33%  "
$ws.Range("D10").Value = "36.99"
$ws.Range("E10").Value = "  +12.68%  "
$ws.Range("E11").Value = "  +3.17%  "
$ws.Range("D12").Value = "7.45"
$ws.Range("E12").Value = "  +7.26%  "
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").Value = "2.659.83"
$ws.Range("E14").Value = "  +3.76%  "
$ws.Range("D15").Value = "2.307.04"
$ws.Range("E15").Value = "  +4.29%  "
$ws.Range("D16").Value = "14.03"
$ws.Range("D17").Value = "0.822"
$ws.Range("E17").Value = "  +5.78%  "
$ws.Range("D18").Value = "46.896.61"
$ws.Range("E18").Value = "  +6.37%  "
$ws.Range("D19").Value = "13.64"
$ws.Range("E19").Value = "  +24.68%  "
$ws.Range("D20").Value = "0.0₃0949"
$ws.Range("E20").Value = "  +4.40%  "
$ws.Range("D21").Value = "6.12"
$ws.Range("E21").Value = "  +3.90%  "
$ws.Range("D22").Value = "66.90"
$ws.Range("E22").Value = "  +3.68%  "
$ws.Range("D23").Value = "249.73"
$ws.Range("E23").Value = "  +5.12%  "
$ws.Range("D24").Value = "2.95"
$ws.Range("E24").Value = "  +5.84%  "
$ws.Range("D25").Value = "1.97"
$ws.Range("E25").Value = "  +6.50%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -1.28%  "
$ws.Range("D27").Value = "44.32"
$ws.Range("E27").Value = "  +16.04%  "
$ws.Range("E28").Value = "  +2.48%  "
$ws.Range("E29").Value = "  +6.84%  "
$ws.Range("D30").Value = "20.24"
$ws.Range("E30").Value = "  +3.98%  "
$ws.Range("D31").Value = "5.81"
$ws.Range("E31").Value = "  +8.22%  "
$ws.Range("E32").Value = "  +8.69%  "
$ws.Range("D33").Value = "148.26"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").Value = "2.69"
$ws.Range("E34").Value = "  +5.83%  "
$ws.Range("D35").Value = "3.20"
$ws.Range("E35").Value = "  +14.35%  "
$ws.Range("D36").Value = "0.114"
$ws.Range("E36").Value = "  +11.63%  "
$ws.Range("E37").Value = "  +3.50%  "
$ws.Range("E38").Value = "  +8.29%  "
$ws.Range("D39").Value = "16.20"
$ws.Range("E39").Value = "  +24.00%  "
$ws.Range("D40").Value = "4.04"
$ws.Range("E40").Value = "  +15.49%  "
$ws.Range("E41").Value = "  +9.87%  "
$ws.Range("E42").Value = "  +1.91%  "
$ws.Range("E43").Value = "  +14.15%  "
$ws.Range("D44").Value = "0.999"
$ws.Range("D45").Value = "1.862.92"
$ws.Range("E45").Value = "  +2.82%  "
$ws.Range("D46").Value = "89.04"
$ws.Range("E46").Value = "  +20.07%  "
$ws.Range("D47").Value = "0.197"
$ws.Range("E47").Value = "  +10.51%  "
$ws.Range("D48").Value = "74.94"
$ws.Range("E48").Value = "  +12.88%  "
$ws.Range("E49").Value = "  +11.67%  "
$ws.Range("D50").Value = "97.35"
$ws.Range("E50").Value = "  +3.83%  "
$ws.Range("D51").Value = "8.07"
$ws.Range("E51").Value = "  +5.89%  "

# Restore original (default) style on column D so the workbook-level
# style table matches the source edit (no leftover custom style index).
$ws.Range("D2:D51").Style = "Normal"
